$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (reflects "Through 2022-02-10" -> "Through 2022-02-11")
$ws.Name = "Through 2022-02-11"

# Update the shared string text for the February label in column A, row 3
$ws.Range("A3").Value = "February (through 02-11)"

# Update cell I2
$ws.Range("I2").Value = 161

# Update row 3 (B3:I3)
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = 28
$ws.Range("E3").Value = 22
$ws.Range("F3").Value = 10
$ws.Range("G3").Value = 26
$ws.Range("H3").Value = 58
$ws.Range("I3").Value = 52

# Update row 4 (B4:I4)
$ws.Range("B4").Value = 30
$ws.Range("C4").Value = 61
$ws.Range("D4").Value = 103
$ws.Range("E4").Value = 108
$ws.Range("F4").Value = 59
$ws.Range("G4").Value = 100
$ws.Range("H4").Value = 275
$ws.Range("I4").Value = 213
